$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Active cases"

$ws.Cells.Item(2, 1).Value = "21st Birthday Party 27 Nov Middels Drouin"
$ws.Cells.Item(2, 2).Value = 10
$ws.Cells.Item(3, 1).Value = "3323 Villa Maria Catholic Homes St Bernadette's Aged Care Sunshine North"
$ws.Cells.Item(3, 2).Value = 12
$ws.Cells.Item(4, 1).Value = "3398 BlueCross Elly Kay Mordialloc"
$ws.Cells.Item(4, 2).Value = 41
$ws.Cells.Item(5, 1).Value = "3601 Baptcare Westhaven community"
$ws.Cells.Item(5, 2).Value = 17
$ws.Cells.Item(6, 1).Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Cells.Item(6, 2).Value = 23
$ws.Cells.Item(7, 1).Value = "3975 Aurrum Aged Care Brunswick West"
$ws.Cells.Item(7, 2).Value = 12
$ws.Cells.Item(8, 1).Value = "4257 BlueCross The Gables Camberwell"
$ws.Cells.Item(8, 2).Value = 20
$ws.Cells.Item(9, 1).Value = "4295 Hope Aged Care Sunshine West"
$ws.Cells.Item(9, 2).Value = 30
$ws.Cells.Item(10, 1).Value = "44087 Fitzroy Primary School Fitzroy"
$ws.Cells.Item(10, 2).Value = 11
$ws.Cells.Item(11, 1).Value = "44444 Nar Nar Goon Primary School Nar Nar Goon"
$ws.Cells.Item(11, 2).Value = 23
$ws.Cells.Item(12, 1).Value = "44666 Gardenvale Primary School Senior School Campus Brighton Eas"
$ws.Cells.Item(12, 2).Value = 18
$ws.Cells.Item(13, 1).Value = "44811 Dandenong North Primary School Dandenong"
$ws.Cells.Item(13, 2).Value = 20
$ws.Cells.Item(14, 1).Value = "44950 Templestowe Valley Primary School Templestowe Lower"
$ws.Cells.Item(14, 2).Value = 71
$ws.Cells.Item(15, 1).Value = "44982 Diamond Creek East Primary School Diamond Creek"
$ws.Cells.Item(15, 2).Value = 16
$ws.Cells.Item(16, 1).Value = "45026 Churchill North Primary School Churchill"
$ws.Cells.Item(16, 2).Value = 13
$ws.Cells.Item(17, 1).Value = "45248 Brookside P-9 College Caroline Springs"
$ws.Cells.Item(17, 2).Value = 23
$ws.Cells.Item(18, 1).Value = "45267 Epping Views Primary School Epping"
$ws.Cells.Item(18, 2).Value = 14
$ws.Cells.Item(19, 1).Value = "45315 Red Hill Consolidated School Red Hill"
$ws.Cells.Item(19, 2).Value = 12
$ws.Cells.Item(20, 1).Value = "45585 Mount Ridley College Craigieburn"
$ws.Cells.Item(20, 2).Value = 14
$ws.Cells.Item(21, 1).Value = "45648 St Brendans Primary School Shepparton"
$ws.Cells.Item(21, 2).Value = 28
$ws.Cells.Item(22, 1).Value = "4574 Village Glen Aged Care Residences Mornington"
$ws.Cells.Item(22, 2).Value = 11
$ws.Cells.Item(23, 1).Value = "45755 St Patricks Catholic Parish Primary School Mentone"
$ws.Cells.Item(23, 2).Value = 17
$ws.Cells.Item(24, 1).Value = "45846 St Mary's School Mooroopna"
$ws.Cells.Item(24, 2).Value = 19
$ws.Cells.Item(25, 1).Value = "45950 St Luke's Primary School Lalor"
$ws.Cells.Item(25, 2).Value = 20
$ws.Cells.Item(26, 1).Value = "46052 St. Francis of Assisi Primary School Mill Park"
$ws.Cells.Item(26, 2).Value = 22
$ws.Cells.Item(27, 1).Value = "46105 Christ the Priest Primary School Caroline Springs"
$ws.Cells.Item(27, 2).Value = 48
$ws.Cells.Item(28, 1).Value = "46115 St Luke's Catholic Primary School Shepparton North"
$ws.Cells.Item(28, 2).Value = 11
$ws.Cells.Item(29, 1).Value = "46117 Marymede Catholic College South Morang"
$ws.Cells.Item(29, 2).Value = 11
$ws.Cells.Item(30, 1).Value = "46221 Bialik College Hawthorn"
$ws.Cells.Item(30, 2).Value = 12
$ws.Cells.Item(31, 1).Value = "46239 Gilson College Taylors Hill"
$ws.Cells.Item(31, 2).Value = 12
$ws.Cells.Item(32, 1).Value = "46287 Oakleigh Grammar Melbourne Private School Oakleigh"
$ws.Cells.Item(32, 2).Value = 11
$ws.Cells.Item(33, 1).Value = "46390 Al Siraat College Epping"
$ws.Cells.Item(33, 2).Value = 11
$ws.Cells.Item(34, 1).Value = "Covenant College Bell Post Hill"
$ws.Cells.Item(34, 2).Value = 17
$ws.Cells.Item(35, 1).Value = "House Party 27 November Private Residence Brunswick West"
$ws.Cells.Item(35, 2).Value = 25
$ws.Cells.Item(36, 1).Value = "Islamic College of Melbourne Tarneit Oct Nov"
$ws.Cells.Item(36, 2).Value = 12
$ws.Cells.Item(37, 1).Value = "Springside Primary School Caroline Springs Nov"
$ws.Cells.Item(37, 2).Value = 33
$ws.Cells.Item(38, 1).Value = "St Vincents Hospital Melbourne Emergency Department Fitzroy"
$ws.Cells.Item(38, 2).Value = 10
$ws.Cells.Item(39, 1).Value = "The Village Early Learning Centre Sandringham"
$ws.Cells.Item(39, 2).Value = 16
$ws.Cells.Item(40, 1).Value = "Torquay Hotel Torquay"
$ws.Cells.Item(40, 2).Value = 14
$ws.Cells.Item(41, 1).Value = "Werribee Mercy Hospital Emergency Department"
$ws.Cells.Item(41, 2).Value = 10

$ws.Range("A42:B69").EntireRow.Delete()
